$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update the A321neo Basic Operating Weight (BOW) value, which drives the
# recalculation of the D20/E20/F20 moment-of-inertia formulas.
$ws.Range("B23").Value = 114640

# Update the active selection to match the new state (F20 instead of F21).
$ws.Range("F20").Select()
